$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.577899385560676
$ws.Range("C2").Value = 0.574479461533623
$ws.Range("D2").Value = 0.652481803034289
$ws.Range("E2").Value = 0.588805087252292
$ws.Range("F2").Value = 0.434701727834531
